# "PostFixProgram Condense into one class"
#
# Para 2: "The code features three classes: A Calculator, a Parser and the
#          main class PostFixProgram. This was done ... aesthetically pleasing."
#      -> "The code features three parts: the Calculator, the Parser and the
#          main class PostFixProgram. "
#
# Para 3: "The CSVParser class is the first that is used, in order to read..."
#      -> "The parser is the first part used, to read..."

$d = $word.ActiveDocument

# 1) "classes: A Calculator, a " -> "parts: the Calculator, the " (leaves the
#    existing "Parser" run / its gramStart-gramEnd proofErr pair untouched).
$d.Content.Find.Execute(
    "The code features three classes: A Calculator, a ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The code features three parts: the Calculator, the ",
    2)

# 2) Drop the long justification sentence (including the second
#    "PostFixProgram" occurrence and its spellStart/spellEnd proofErr pair),
#    keeping just the closing ". ".
$d.Content.Find.Execute(
    ". This was done just for the sake of not making the main class look too bloated, it would have been just as easy to keep the two methods within the PostFixProgram, but I felt it just was a tad more aesthetically pleasing.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ". ",
    2)

# 3) "The CSVParser class is the first that is used, in order to " ->
#    "The parser is the first part used, to " (removes the spellStart/spellEnd
#    around "CSVParser" and the gramStart/gramEnd around "in order to").
$d.Content.Find.Execute(
    "The CSVParser class is the first that is used, in order to ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The parser is the first part used, to ",
    2)
